# --- Update biz plan, biz model ---
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Content (order matches shared-string build order of the original edit) --
$ws.Range("C1").Value = "This document indicates the roles of each member in the team. The role can be changed afterwards to fit with the condition of the company"

$ws.Range("A2").Value = "Overview"
$ws.Range("B2").Value = "Geographic roles"

$ws.Range("A3").Value = "CEO"

$ws.Range("A4").Value = "Functional roles include:"
$ws.Range("A5").Value = "Area manager"
$ws.Range("A6").Value = "Finance"
$ws.Range("A7").Value = "Marketing"
$ws.Range("A8").Value = "Sales"
$ws.Range("A9").Value = "Negotiation"
$ws.Range("A10").Value = "Production"
$ws.Range("A11").Value = "Logistics"
$ws.Range("A12").Value = "R&D"
$ws.Range("A13").Value = "Etc."

$ws.Range("B3").Value = "Home Office (Leich.)"
$ws.Range("C3").Value = "US"
$ws.Range("D3").Value = "EU"
$ws.Range("E3").Value = "Brazil"

# Formatting ---------------------------------------------------------------
# 16pt heading font on every cell that actually holds content (this becomes
# cellXfs index 1, reused by every plain cell in the sheet).
$ws.Range("C1").Font.Size = 16
$ws.Range("A2").Font.Size = 16
$ws.Range("A3:E3").Font.Size = 16
$ws.Range("A4").Font.Size = 16
$ws.Range("A5").Font.Size = 16
$ws.Range("A6").Font.Size = 16
$ws.Range("A7").Font.Size = 16
$ws.Range("A8").Font.Size = 16
$ws.Range("A9").Font.Size = 16
$ws.Range("A10").Font.Size = 16
$ws.Range("A11").Font.Size = 16
$ws.Range("A12").Font.Size = 16
$ws.Range("A13").Font.Size = 16

# B2:E2 get the 16pt font too, then centered and merged - becomes cellXfs index 2.
$ws.Range("B2:E2").Font.Size = 16
$ws.Range("B2:E2").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B2:E2").Merge()

# Column widths (best-fit-ish, closest achievable through this host) -------
$ws.Columns.Item(1).ColumnWidth = 27.25
$ws.Columns.Item(2).ColumnWidth = 22.75

$ws.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------------
# The source workbook also carries the (unused, built-in) "Hyperlink" /
# "Followed Hyperlink" cell styles in its style table - a leftover of the
# Mac Excel template this was saved from. Materialize + immediately undo them
# on a scratch cell so the style/font entries land in the style table without
# actually formatting any real content.
$scratch1 = $ws.Range("Z100")
$scratch1.Value = "tmp"
$scratch1.Style = "Hyperlink"
$excel.Undo() | Out-Null

$scratch2 = $ws.Range("Z100")
$scratch2.Value = "tmp"
$scratch2.Style = "Followed Hyperlink"
$excel.Undo() | Out-Null
